$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2021 column of data, mirroring the styles of the neighboring 2020 column (Q)
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

$ws.Range("Q5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 72

$ws.Application.CutCopyMode = $false

# Update the view: clear the frozen/top-left cell and reset selection to R1
$ws.Range("R1").Select()
